# Config update for Edge Browser:
# Mark the "Edge" column (C) as "Yes" (was "No") for the applicable test
# rows on the "Automation Tests" sheet, and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,18,19,20,21,22,23,24,25,26,27,29,30,32,33,34,35,36)
foreach ($r in $rows) {
    $ws.Range("C$r").Value = "Yes"
}

# Move the active selection to B27 (and drop the previous frozen/scrolled
# top-left cell, matching the saved view state).
$ws.Range("B27").Select()
